$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Batch 1 (rows 2-11, Iterations = 100) ---
$ws.Range("D2").Value = 0.06631325354785259
$ws.Range("E2").Value = 0.06631325354785259

$ws.Range("D3").Value = 0.0005176182402187505
$ws.Range("E3").Value = 0.0005176182402187505

$ws.Range("D4").Value = 0.9767506921254677
$ws.Range("E4").Value = 0.9767506921254677

$ws.Range("D5").Value = [double]"5.753766652052312E-05"
$ws.Range("E5").Value = [double]"5.753766652052312E-05"

$ws.Range("D6").Value = 0.0395370218682276
$ws.Range("E6").Value = 0.0395370218682276

$ws.Range("D7").Value = 0.9773565464416615
$ws.Range("E7").Value = 0.02264345355833852

$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"1.75047130808624E-08"
$ws.Range("E8").Value = 0.9999999824952869

$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"1.02578954508364E-14"
$ws.Range("E9").Value = 0.9999999999999898

$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"8.083883234871182E-08"
$ws.Range("E10").Value = 0.9999999191611677

$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.07603773225776467
$ws.Range("E11").Value = 0.9239622677422353
$ws.Range("F11").Value = 7.287277221679688
$ws.Range("G11").Value = 0.5

# --- Batch 2 (rows 12-21, Iterations = 200) ---
$ws.Range("D12").Value = 0.05291952315915446
$ws.Range("E12").Value = 0.05291952315915446

$ws.Range("D13").Value = 0.0008039404600372358
$ws.Range("E13").Value = 0.0008039404600372358

$ws.Range("D14").Value = 0.9987484147992824
$ws.Range("E14").Value = 0.9987484147992824

$ws.Range("D15").Value = [double]"1.016026132827289E-05"
$ws.Range("E15").Value = [double]"1.016026132827289E-05"

$ws.Range("D16").Value = 0.009859497460229901
$ws.Range("E16").Value = 0.009859497460229901

$ws.Range("D17").Value = 0.9956815709578656
$ws.Range("E17").Value = 0.00431842904213442

$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"3.026116525208176E-10"
$ws.Range("E18").Value = 0.9999999996973884

$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"4.963012513862474E-20"
$ws.Range("E19").Value = 1

$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"1.144764449829875E-09"
$ws.Range("E20").Value = 0.9999999988552356

$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.005976142633226013
$ws.Range("E21").Value = 0.994023857366774
$ws.Range("F21").Value = 9.882906913757324
$ws.Range("G21").Value = 0.5
